$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): A1 tab title text changes ---
$ws.Range("A1").Value = "TC02_Bento_LocalSearch-UploadCaseSet_Enter_CASEID_List"

# --- Row 2: CasesTab ---
$casesTabB = @"
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss:study_subject)
	WHERE ss.study_subject_id = 'BENTO-CASE-3405467'
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
RETURN DISTINCT 
	ss.study_subject_id AS ``Case ID``,
	p.program_acronym AS ``Program Code``,
	p.program_id AS ``Program ID``,
	s.study_acronym AS ``Arm``,
	ss.disease_subtype AS ``Diagnosis``,
	sf.grouped_recurrence_score AS ``Recurrence Score``,
	d.tumor_size_group AS ``Tumor Size (cm)``,
	d.er_status AS ``ER Status``,
	d.pr_status AS ``PR Status``,
	demo.age_at_index AS ``Age (years)``,
	demo.survival_time AS ``Survival (days)``
"@

$casesTabC = @"
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)-[:study_subject_of_study]->(s)
WHERE ss.study_subject_id = 'BENTO-CASE-3405467'
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files

"@

$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $casesTabB
$ws.Range("B2").WrapText = $true
$ws.Range("C2").Value = $casesTabC
$ws.Range("C2").WrapText = $true
$ws.Range("D2").Value = "TC02_Bento_LocalSearch-UploadCaseSet_Enter_CASEID_List_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC02_Bento_LocalSearch-UploadCaseSet_Enter_CASEID_List_WebData.xlsx"
$ws.Rows.Item(2).RowHeight = 345

# --- Row 3: SamplesTab (new row) ---
$samplesTabB = @"
MATCH (ss:study_subject)
	WHERE ss.study_subject_id = 'BENTO-CASE-3405467'
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[*..2]-(parent)<--(f:file)
OPTIONAL MATCH (f)-[:file_of_laboratory_procedure]->(lp)
RETURN DISTINCT 
	samp.sample_id AS ``Sample ID``,
	ss.study_subject_id AS ``Case ID``,
	p.program_acronym AS ``Program Code``,
	s.study_acronym AS ``Arm``,
	ss.disease_subtype AS ``Diagnosis``,
	samp.tissue_type AS ``Tissue Type``,
	samp.composition AS ``Tissue Composition``,
	samp.sample_anatomic_site AS ``Sample Anatomic Site``,
	samp.method_of_sample_procurement AS ``Sample Procurement Method``,
	lp.test_name AS ``platform``
"@

$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $samplesTabB
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = $casesTabC
$ws.Range("C3").WrapText = $true
$ws.Range("D3").Value = "TC02_Bento_LocalSearch-UploadCaseSet_Enter_CASEID_List_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC02_Bento_LocalSearch-UploadCaseSet_Enter_CASEID_List_WebData.xlsx"
$ws.Rows.Item(3).RowHeight = 255

# --- Row 4: FilesTab (new row) ---
$filesTabB = @"
MATCH (ss:study_subject)
	WHERE ss.study_subject_id = 'BENTO-CASE-3405467'
MATCH (ss)<-[*..2]-(parent)<--(f:file)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
RETURN DISTINCT 
	f.file_name AS ``File Name``,
	head(labels(parent)) AS ``Association``,
	f.file_description AS ``Description``,
	f.file_format AS ``File Format``,
	f.file_size AS ``Size``,
	p.program_acronym AS ``Program Code``,
	s.study_acronym AS ``Arm``,
	ss.study_subject_id AS ``Case ID``,
	samp.sample_id AS ``Sample ID``
"@

$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $filesTabB
$ws.Range("B4").WrapText = $true
$ws.Range("C4").Value = $casesTabC
$ws.Range("C4").WrapText = $true
$ws.Range("D4").Value = "TC02_Bento_LocalSearch-UploadCaseSet_Enter_CASEID_List_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC02_Bento_LocalSearch-UploadCaseSet_Enter_CASEID_List_WebData.xlsx"
$ws.Rows.Item(4).RowHeight = 255

# --- Selection / view state ---
$ws.Range("D4:E4").Select()